# Replace every occurrence of the (multi-run) "Datas das campanhas de 2018
# que usam Perseu: ..." paragraph text with a single, plain run containing
# "Datas das campanhas de Leo: 14-23 de abril, 14-23 de maio" (no run
# formatting at all, matching the target OOXML).
#
# The original sentence is split across 4 runs (different rPr/rsid on
# each), so a plain Find/Replace would just merge them into one run while
# keeping the first run's rPr. Instead, for each match we clear the whole
# found range's text (which removes the runs entirely) and then insert
# the new text after it, which creates a brand-new run with no rPr.

$d = $word.ActiveDocument

$oldText = "Datas das campanhas de 2018 que usam Perseu: 30 de outubro a 8 de novembro e 29 de novembro a 8 de dezembro"
$newText = "Datas das campanhas de Leo: 14-23 de abril, 14-23 de maio"

$replacements = 0
$maxIter = 50

while ($replacements -lt $maxIter) {
    $r = $d.Content
    $found = $r.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        break
    }

    $r.Text = ""
    $r.InsertAfter($newText) | Out-Null

    $replacements = $replacements + 1
}

Write-Output "Replaced $replacements occurrence(s) of the campaign-dates paragraph."
